$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-27 07:48:21"
$ws.Range("K2").Value = "0.1 MJ/m2"
$ws.Range("M2").Value = "3.0 °C 7:29 TU"
$ws.Range("E3").Value = "2026-02-27 07:48:23"
$ws.Range("K3").Value = "0.0 MJ/m2"
$ws.Range("E4").Value = "2026-02-27 07:48:26"
$ws.Range("L4").Value = "11.5 km/h - 317º 7:03 TU"
$ws.Range("E5").Value = "2026-02-27 07:48:28"
$ws.Range("H5").Formula = '="37%"'
$ws.Range("H5").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4163) | Out-Null
$ws.Range("K5").Value = "0.0 MJ/m2"
$ws.Range("E6").Value = "2026-02-27 07:48:31"
$ws.Range("L6").Value = "6.8 km/h - 6º 7:06 TU"
$ws.Range("E7").Value = "2026-02-27 07:48:33"
$ws.Range("K7").Value = "0.1 MJ/m2"
$ws.Range("O7").Value = "9.7 °C"
$ws.Range("E8").Value = "2026-02-27 07:48:35"
$ws.Range("K8").Value = "0.2 MJ/m2"
$ws.Range("E9").Value = "2026-02-27 07:48:38"
$ws.Range("K9").Value = "0.1 MJ/m2"
$ws.Range("M9").Value = "9.9 °C 7:29 TU"
$ws.Range("O9").Value = "8.7 °C"
$ws.Range("E10").Value = "2026-02-27 07:48:39"
$ws.Range("K10").Value = "0.1 MJ/m2"
$ws.Range("O10").Value = "8.8 °C"
$ws.Range("E11").Value = "2026-02-27 07:48:40"
$ws.Range("E12").Value = "2026-02-27 07:48:41"
$ws.Range("O12").Value = "8.6 °C"
$ws.Range("E13").Value = "2026-02-27 07:48:42"
$ws.Range("J13").Value = "1032.4 hPa"
$ws.Range("K13").Value = "0.1 MJ/m2"
$ws.Range("O13").Value = "-1.9 °C"
$ws.Range("E14").Value = "2026-02-27 07:48:43"
$ws.Range("K14").Value = "0.1 MJ/m2"
$ws.Range("E15").Value = "2026-02-27 07:48:44"
$ws.Range("M15").Value = "9.9 °C 7:20 TU"
$ws.Range("O15").Value = "8.8 °C"
$ws.Range("E16").Value = "2026-02-27 07:48:45"
$ws.Range("G16").Value = "67 cm"
$ws.Range("K16").Value = "0.0 MJ/m2"
$ws.Range("E17").Value = "2026-02-27 07:48:46"
$ws.Range("K17").Value = "0.2 MJ/m2"
$ws.Range("E18").Value = "2026-02-27 07:48:48"
$ws.Range("H18").Formula = '="95%"'
$ws.Range("H18").Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "2026-02-27 07:48:49"
$ws.Range("K19").Value = "0.0 MJ/m2"
$ws.Range("E20").Value = "2026-02-27 07:48:50"
$ws.Range("H20").Formula = '="56%"'
$ws.Range("H20").Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4163) | Out-Null
$ws.Range("K20").Value = "0.1 MJ/m2"
$ws.Range("M20").Value = "3.6 °C 7:25 TU"
$ws.Range("E21").Value = "2026-02-27 07:48:51"
$ws.Range("J21").Value = "1029.2 hPa"
$ws.Range("K21").Value = "0.1 MJ/m2"
$ws.Range("O21").Value = "3.4 °C"
$ws.Range("E22").Value = "2026-02-27 07:48:53"
$ws.Range("K22").Value = "0.1 MJ/m2"
$ws.Range("O22").Value = "0.6 °C"
$ws.Range("E23").Value = "2026-02-27 07:48:56"
$ws.Range("H23").Formula = '="39%"'
$ws.Range("H23").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4163) | Out-Null
$ws.Range("K23").Value = "0.0 MJ/m2"
$ws.Range("E24").Value = "2026-02-27 07:48:58"
$ws.Range("K24").Value = "0.1 MJ/m2"
$ws.Range("O24").Value = "3.7 °C"
$ws.Range("E25").Value = "2026-02-27 07:49:00"
$ws.Range("K25").Value = "0.2 MJ/m2"
$ws.Range("M25").Value = "6.7 °C 7:12 TU"
$ws.Range("O25").Value = "4.7 °C"
$ws.Range("E26").Value = "2026-02-27 07:49:03"
$ws.Range("H26").Formula = '="42%"'
$ws.Range("H26").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4163) | Out-Null
$ws.Range("K26").Value = "0.0 MJ/m2"
$ws.Range("M26").Value = "9.8 °C 7:29 TU"
$ws.Range("E27").Value = "2026-02-27 07:49:05"
$ws.Range("H27").Formula = '="42%"'
$ws.Range("H27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4163) | Out-Null
$ws.Range("K27").Value = "0.1 MJ/m2"
$ws.Range("M27").Value = "5.3 °C 7:04 TU"
$ws.Range("E28").Value = "2026-02-27 07:49:08"
$ws.Range("J28").Value = "1026.2 hPa"
$ws.Range("K28").Value = "0.1 MJ/m2"
$ws.Range("L28").Value = "11.2 km/h - 281º 7:12 TU"
$ws.Range("E29").Value = "2026-02-27 07:49:10"
$ws.Range("K29").Value = "0.1 MJ/m2"
$ws.Range("L29").Value = "10.4 km/h - 355º 7:19 TU"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-27 07:49:13"
$ws.Range("K30").Value = "0.1 MJ/m2"
$ws.Range("E31").Value = "2026-02-27 07:49:15"
$ws.Range("I31").Value = "0.1 mm"
$ws.Range("E32").Value = "2026-02-27 07:49:17"
$ws.Range("K32").Value = "0.1 MJ/m2"
$ws.Range("M32").Value = "2.9 °C 7:29 TU"
$ws.Range("E33").Value = "2026-02-27 07:49:19"
$ws.Range("K33").Value = "0.1 MJ/m2"
$ws.Range("N33").Value = "0.3 °C 7:06 TU"
$ws.Range("E34").Value = "2026-02-27 07:49:22"
$ws.Range("H34").Formula = '="43%"'
$ws.Range("H34").Copy() | Out-Null
$ws.Range("H34").PasteSpecial(-4163) | Out-Null
$ws.Range("L34").Value = "21.6 km/h - 22º 7:19 TU"
$ws.Range("M34").Value = "4.6 °C 7:13 TU"
$ws.Range("O34").Value = "2.4 °C"
$ws.Range("E35").Value = "2026-02-27 07:49:24"
$ws.Range("K35").Value = "0.0 MJ/m2"
$ws.Range("O35").Value = "9.8 °C"
$ws.Range("E36").Value = "2026-02-27 07:49:27"
$ws.Range("M36").Value = "10.9 °C 7:29 TU"
$ws.Range("E37").Value = "2026-02-27 07:49:29"
$ws.Range("H37").Formula = '="90%"'
$ws.Range("H37").Copy() | Out-Null
$ws.Range("H37").PasteSpecial(-4163) | Out-Null
$ws.Range("J37").Value = "1028.6 hPa"
$ws.Range("N37").Value = "0.0 °C 7:00 TU"
$ws.Range("O37").Value = "2.2 °C"
$ws.Range("E38").Value = "2026-02-27 07:49:32"
$ws.Range("E39").Value = "2026-02-27 07:49:34"
$ws.Range("K39").Value = "0.2 MJ/m2"
$ws.Range("E40").Value = "2026-02-27 07:49:36"
$ws.Range("J40").Value = "1029.9 hPa"
$ws.Range("N40").Value = "0.4 °C 7:05 TU"
$ws.Range("E41").Value = "2026-02-27 07:49:39"
$ws.Range("K41").Value = "0.1 MJ/m2"
$ws.Range("O41").Value = "8.2 °C"
$ws.Range("E42").Value = "2026-02-27 07:49:41"
$ws.Range("M42").Value = "10.0 °C 7:29 TU"
$ws.Range("E43").Value = "2026-02-27 07:49:44"
$ws.Range("K43").Value = "0.2 MJ/m2"
$ws.Range("O43").Value = "3.5 °C"
$ws.Range("E44").Value = "2026-02-27 07:49:46"
$ws.Range("H44").Formula = '="66%"'
$ws.Range("H44").Copy() | Out-Null
$ws.Range("H44").PasteSpecial(-4163) | Out-Null
$ws.Range("K44").Value = "0.0 MJ/m2"
$ws.Range("N44").Value = "-2.2 °C 7:19 TU"
$ws.Range("O44").Value = "-0.5 °C"
$ws.Range("E45").Value = "2026-02-27 07:49:49"
$ws.Range("J45").Value = "1026.3 hPa"
$ws.Range("K45").Value = "0.0 MJ/m2"
$ws.Range("E46").Value = "2026-02-27 07:49:51"
$ws.Range("J46").Value = "1026.0 hPa"
$ws.Range("K46").Value = "0.0 MJ/m2"
$ws.Range("O46").Value = "6.1 °C"
$excel.CutCopyMode = $false
